# PROD-10305: add/delete row features
# Reflects the post-edit state of cypress/fixtures/customers.xlsx:
#  - Howard's birthday (E5) corrected from 21/05/2002 to 12/05/1987
#  - Billy's canDrinkAlcohol (C6) stored as a TRUE() formula instead of a
#    literal boolean
#  - active selection on the "customers" sheet left on E5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customers")

# Correct Howard's birthday
$ws.Range("E5").Value = "12/05/1987"

# Store Billy's canDrinkAlcohol value as a formula (=TRUE()) rather than a
# bare boolean literal
$ws.Range("C6").Formula = "=TRUE()"

# Leave the active selection on E5
$ws.Range("E5").Select()
